$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh snapshot.
# D-column numeric-looking values get a leading apostrophe so Excel
# keeps them as text (matching the sheet's existing text-stored numbers)
# instead of auto-converting them to the Number type.

# D2: "63.824.07" -> "63.838.08"; E2: "  -0.08%  " -> "  -0.05%  "
$ws.Range("D2").Value = "63.838.08"
$ws.Range("E2").Value = "  -0.05%  "

# D3: "2.734.00" -> "2.733.66"; E3: "  -0.61%  " -> "  -0.53%  "
$ws.Range("D3").Value = "2.733.66"
$ws.Range("E3").Value = "  -0.53%  "

# E4: "  +0.11%  " -> "  -0.01%  "
$ws.Range("E4").Value = "  -0.01%  "

# D5: "565.28" -> "565.10"; E5: "  -1.68%  " -> "  -1.57%  "
$ws.Range("D5").Value = "'565.10"
$ws.Range("E5").Value = "  -1.57%  "

# D6: "160.53" -> "160.46"; E6: "  +1.83%  " -> "  +1.93%  "
$ws.Range("D6").Value = "'160.46"
$ws.Range("E6").Value = "  +1.93%  "

# E7: "  +0.01%  " -> "  -0.04%  "
$ws.Range("E7").Value = "  -0.04%  "

# E8: "  -0.95%  " -> "  -0.89%  "
$ws.Range("E8").Value = "  -0.89%  "

# E9: "  -0.14%  " -> "  -0.04%  "
$ws.Range("E9").Value = "  -0.04%  "

# E10: "  +4.10%  " -> "  +3.99%  "
$ws.Range("E10").Value = "  +3.99%  "

# D11: "5.62" -> "5.60"; E11: "  -1.69%  " -> "  -1.03%  "
$ws.Range("D11").Value = "'5.60"
$ws.Range("E11").Value = "  -1.03%  "

# E12: "  -1.49%  " -> "  -1.47%  "
$ws.Range("E12").Value = "  -1.47%  "

# D13: "3.219.37" -> "3.218.36"; E13: "  -0.54%  " -> "  -0.58%  "
$ws.Range("D13").Value = "3.218.36"
$ws.Range("E13").Value = "  -0.58%  "

# D14: "26.91" -> "26.90"; E14: "  +1.51%  " -> "  +1.21%  "
$ws.Range("D14").Value = "'26.90"
$ws.Range("E14").Value = "  +1.21%  "

# D15: "63.681.29" -> "63.664.27"; E15: "  +0.28%  " -> "  +0.23%  "
$ws.Range("D15").Value = "63.664.27"
$ws.Range("E15").Value = "  +0.23%  "

# E16: "  -0.84%  " -> "  -0.87%  "
$ws.Range("E16").Value = "  -0.87%  "

# D17: "2.743.35" -> "2.738.36"; E17: "  -0.36%  " -> "  -0.53%  "
$ws.Range("D17").Value = "2.738.36"
$ws.Range("E17").Value = "  -0.53%  "

# D18: "12.33" -> "12.35"; E18: "  +1.96%  " -> "  +2.12%  "
$ws.Range("D18").Value = "'12.35"
$ws.Range("E18").Value = "  +2.12%  "

# D20: "355.78" -> "355.99"; E20: "  +0.09%  " -> "  +0.22%  "
$ws.Range("D20").Value = "'355.99"
$ws.Range("E20").Value = "  +0.22%  "

# D21: "6.61" -> "6.60"; E21: "  -1.26%  " -> "  -1.40%  "
$ws.Range("D21").Value = "'6.60"
$ws.Range("E21").Value = "  -1.40%  "

# E22: "  -0.12%  " -> "  -0.18%  "
$ws.Range("E22").Value = "  -0.18%  "

# D23: "0.519" -> "0.520"; E23: "  -2.59%  " -> "  -2.57%  "
$ws.Range("D23").Value = "'0.520"
$ws.Range("E23").Value = "  -2.57%  "

# D24: "64.17" -> "64.14"; E24: "  -1.51%  " -> "  -1.56%  "
$ws.Range("D24").Value = "'64.14"
$ws.Range("E24").Value = "  -1.56%  "

# E25: "  -0.30%  " -> "  -0.39%  "
$ws.Range("E25").Value = "  -0.39%  "

# D26: "1.00" -> "0.999"; E26: "  +0.10%  " -> "  +0.07%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.07%  "

# D27: "8.36" -> "8.35"; E27: "  -1.11%  " -> "  -1.25%  "
$ws.Range("D27").Value = "'8.35"
$ws.Range("E27").Value = "  -1.25%  "

# E28: "  +0.13%  " -> "  +0.23%  "
$ws.Range("E28").Value = "  +0.23%  "

# E29: "  +2.46%  " -> "  +2.47%  "
$ws.Range("E29").Value = "  +2.47%  "

# E30: "  +11.21%  " -> "  +11.01%  "
$ws.Range("E30").Value = "  +11.01%  "

# D31: "7.17" -> "7.18"; E31: "  +1.44%  " -> "  +1.56%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  +1.56%  "

# D32: "167.41" -> "167.34"; E32: "  -0.55%  " -> "  -0.60%  "
$ws.Range("D32").Value = "'167.34"
$ws.Range("E32").Value = "  -0.60%  "

# D33: "4.91" -> "4.92"; E33: "  -0.02%  " -> "  +0.26%  "
$ws.Range("D33").Value = "'4.92"
$ws.Range("E33").Value = "  +0.26%  "

# E34: "  +2.23%  " -> "  +2.33%  "
$ws.Range("E34").Value = "  +2.33%  "

# D35: "20.04" -> "20.06"; E35: "  -0.51%  " -> "  -0.50%  "
$ws.Range("D35").Value = "'20.06"
$ws.Range("E35").Value = "  -0.50%  "

# E37: "  +1.06%  " -> "  +1.18%  "
$ws.Range("E37").Value = "  +1.18%  "

# D38: "0.978" -> "0.976"; E38: "  -0.90%  " -> "  -0.36%  "
$ws.Range("D38").Value = "'0.976"
$ws.Range("E38").Value = "  -0.36%  "

# D39: "347.46" -> "346.61"; E39: "  +4.96%  " -> "  +4.36%  "
$ws.Range("D39").Value = "'346.61"
$ws.Range("E39").Value = "  +4.36%  "

# D40: "6.29" -> "6.30"; E40: "  +2.30%  " -> "  +2.41%  "
$ws.Range("D40").Value = "'6.30"
$ws.Range("E40").Value = "  +2.41%  "

# E41: "  -1.80%  " -> "  -1.71%  "
$ws.Range("E41").Value = "  -1.71%  "

# D42: "38.61" -> "38.65"; E42: "  -0.81%  " -> "  -0.82%  "
$ws.Range("D42").Value = "'38.65"
$ws.Range("E42").Value = "  -0.82%  "

# D43: "21.84" -> "21.81"; E43: "  +1.57%  " -> "  +1.60%  "
$ws.Range("D43").Value = "'21.81"
$ws.Range("E43").Value = "  +1.60%  "

# B44: "Hedera" -> "EnergySwap"; C44: "https://coinranking.com/coin/jad286TjB+hedera-hbar" -> "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D44: "0.0582" -> "20.90"; E44: "  -0.72%  " -> "  -2.64%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'20.90"
$ws.Range("E44").Value = "  -2.64%  "

# B45: "EnergySwap" -> "Hedera"; C45: "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" -> "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D45: "20.88" -> "0.0582"; E45: "  -2.72%  " -> "  -0.79%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0582"
$ws.Range("E45").Value = "  -0.79%  "

# E46: "  +0.94%  " -> "  +1.06%  "
$ws.Range("E46").Value = "  +1.06%  "

# E48: "  -1.39%  " -> "  -1.41%  "
$ws.Range("E48").Value = "  -1.41%  "

# D49: "132.46" -> "132.36"; E49: "  -1.63%  " -> "  -1.66%  "
$ws.Range("D49").Value = "'132.36"
$ws.Range("E49").Value = "  -1.66%  "

# E50: "  -0.05%  " -> "  -0.10%  "
$ws.Range("E50").Value = "  -0.10%  "

# D51: "11.08" -> "11.07"; E51: "  +0.56%  " -> "  +0.47%  "
$ws.Range("D51").Value = "'11.07"
$ws.Range("E51").Value = "  +0.47%  "
